# Fix variablechange to use self instead of local variables
#
# - Solar PV sheet: price rounding F2 179.69 -> 180, F6 177.7 -> 178
# - Solar Thermal sheet: G2 "self" price becomes the numeric 1968 (was the
#   locally-duplicated shared string " 1 968.00"), and G5's redundant
#   `=999` formula becomes a plain literal 999 (still using "self"/value,
#   not a recomputed local formula)
# - Active tab / selection moves from "Solar PV" (B10) to "Wind Turbine" (F7),
#   and the lingering selection left on "Solar PV" becomes F12

$wb = $excel.ActiveWorkbook

$solarPV      = $wb.Worksheets.Item("Solar PV")
$solarThermal = $wb.Worksheets.Item("Solar Thermal")
$windTurbine  = $wb.Worksheets.Item("Wind Turbine")

# --- Data fixes -----------------------------------------------------------

# Solar PV: round the prices
$solarPV.Range("F2").Value = 180
$solarPV.Range("F6").Value = 178

# Solar Thermal: store the price as a real number instead of a separate
# shared-string literal, and drop the pointless "=999" formula in favour of
# the literal value it always evaluated to
$solarThermal.Range("G2").Value = 1968
$solarThermal.Range("G5").Value = 999

# --- Selection / active sheet ---------------------------------------------

# Leave a stale selection on Solar PV (it's no longer the active tab)
[void]$solarPV.Range("F12").Select()

# Wind Turbine becomes the active (selected) tab
[void]$windTurbine.Activate()
[void]$windTurbine.Range("F7").Select()
